# Auto-generated edit script applying cryptos list price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub3 = [char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.839.46'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.334.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.89%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.612'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.329.59'
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.620'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.159'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000270'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.97'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.874.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("E17").Value = '  -2.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.340.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '63.800.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.975'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '432.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.58'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.51'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '591.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.40'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.78%  '
$ws.Range("E34").Value = '  -1.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.143'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.47'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.27'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0' + $sub3.ToString() + '0746'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.364'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.108.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.83'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("E46").Value = '  -0.96%  '
$ws.Range("E47").Value = '  -1.76%  '
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("E49").Value = '  -3.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.17'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '133.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.04%  '

Write-Host "Applied" 88 "cell updates"
